$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously held "pregunta1" (first question) with answers
# resa1/resb1/resc1/resd1. Replace it with a new question about the
# Colombian flag, and drop the now-unused strings from the shared table.
$ws.Range("A2").Value = "rojo"
$ws.Range("B2").Value = "verde"
$ws.Range("C2").Value = "azul"
$ws.Range("D2").Value = "amarillo"
$ws.Range("E2").Value = "primer color del bandera de colombia"
$ws.Range("F2").Value = "d"
$ws.Range("G2").Value = "el color del oro"

# Column width adjustments (best-fit) so the longer question/hint text fits.
$ws.Columns.Item(5).ColumnWidth = 33.8
$ws.Columns.Item(7).ColumnWidth = 13.5

# Update the active selection to reflect where the author left off editing.
$ws.Range("G2").Select()

$wb.Save()
